$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2140575079872205
$ws.Range("C2").Value = 0.549520766773163
$ws.Range("J2").Value = 0.01597444089456869
$ws.Range("P2").Value = 0.1373801916932907
$ws.Range("S2").Value = 0.08306709265175719
$ws.Range("B3").Value = 0.0111731843575419
$ws.Range("C3").Value = 0.01675977653631285
$ws.Range("J3").Value = 0.0223463687150838
$ws.Range("P3").Value = 0.7262569832402235
$ws.Range("S3").Value = 0.223463687150838
$ws.Range("P4").Value = 0.7037037037037037
$ws.Range("S4").Value = 0.2962962962962963
$ws.Range("B6").Value = 0.05741626794258373
$ws.Range("D6").Value = 0.009569377990430622
$ws.Range("F6").Value = 0.03827751196172249
$ws.Range("J6").Value = 0.277511961722488
$ws.Range("O6").Value = 0.01913875598086124
$ws.Range("Q6").Value = 0.1435406698564593
$ws.Range("R6").Value = 0.07655502392344497
$ws.Range("S6").Value = 0.3779904306220095
$ws.Range("B7").Value = 0.115606936416185
$ws.Range("D7").Value = 0.02890173410404624
$ws.Range("F7").Value = 0.06936416184971098
$ws.Range("J7").Value = 0.1271676300578035
$ws.Range("Q7").Value = 0.1329479768786127
$ws.Range("R7").Value = 0.1098265895953757
$ws.Range("S7").Value = 0.4161849710982659
$ws.Range("B8").Value = 0.07610993657505286
$ws.Range("D8").Value = 0.02114164904862579
$ws.Range("F8").Value = 0.06342494714587738
$ws.Range("J8").Value = 0.09936575052854123
$ws.Range("O8").Value = 0.01691331923890063
$ws.Range("Q8").Value = 0.1649048625792812
$ws.Range("R8").Value = 0.09090909090909091
$ws.Range("S8").Value = 0.46723044397463
$ws.Range("B9").Value = 0.1149425287356322
$ws.Range("D9").Value = 0.01915708812260536
$ws.Range("E9").Value = 0.003831417624521073
$ws.Range("F9").Value = 0.05747126436781609
$ws.Range("J9").Value = 0.09961685823754789
$ws.Range("O9").Value = 0.01149425287356322
$ws.Range("Q9").Value = 0.1340996168582375
$ws.Range("R9").Value = 0.0842911877394636
$ws.Range("S9").Value = 0.475095785440613
$ws.Range("B10").Value = 0.1135091926458833
$ws.Range("D10").Value = 0.026378896882494
$ws.Range("E10").Value = 0.003996802557953637
$ws.Range("F10").Value = 0.06554756195043965
$ws.Range("J10").Value = 0.1239008792965628
$ws.Range("O10").Value = 0.007993605115907274
$ws.Range("Q10").Value = 0.193445243804956
$ws.Range("R10").Value = 0.07753796962430055
$ws.Range("S10").Value = 0.3876898481215028
$ws.Range("F11").Value = 0.003759398496240601
$ws.Range("G11").Value = 0.1353383458646616
$ws.Range("J11").Value = 0.1428571428571428
$ws.Range("K11").Value = 0.1729323308270677
$ws.Range("L11").Value = 0.5375939849624061
$ws.Range("S11").Value = 0.007518796992481203
$ws.Range("G12").Value = 0.7549668874172185
$ws.Range("J12").Value = 0.1390728476821192
$ws.Range("K12").Value = 0.01324503311258278
$ws.Range("L12").Value = 0.04635761589403974
$ws.Range("S12").Value = 0.04635761589403974
$ws.Range("G13").Value = 0.7317073170731707
$ws.Range("J13").Value = 0.2195121951219512
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.0182648401826484
$ws.Range("H15").Value = 0.2054794520547945
$ws.Range("I15").Value = 0.0958904109589041
$ws.Range("J15").Value = 0.3105022831050228
$ws.Range("K15").Value = 0.0730593607305936
$ws.Range("M15").Value = 0.0182648401826484
$ws.Range("O15").Value = 0.0867579908675799
$ws.Range("S15").Value = 0.1917808219178082
$ws.Range("F16").Value = 0.01442307692307692
$ws.Range("H16").Value = 0.1634615384615385
$ws.Range("I16").Value = 0.125
$ws.Range("J16").Value = 0.3990384615384616
$ws.Range("K16").Value = 0.09134615384615384
$ws.Range("M16").Value = 0.03365384615384615
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.0576923076923077
$ws.Range("S16").Value = 0.1105769230769231
$ws.Range("F17").Value = 0.007389162561576354
$ws.Range("H17").Value = 0.1773399014778325
$ws.Range("I17").Value = 0.1699507389162561
$ws.Range("J17").Value = 0.3940886699507389
$ws.Range("K17").Value = 0.06403940886699508
$ws.Range("M17").Value = 0.01477832512315271
$ws.Range("O17").Value = 0.04679802955665024
$ws.Range("S17").Value = 0.125615763546798
$ws.Range("F18").Value = 0.01515151515151515
$ws.Range("H18").Value = 0.1818181818181818
$ws.Range("I18").Value = 0.1060606060606061
$ws.Range("J18").Value = 0.4191919191919192
$ws.Range("K18").Value = 0.1212121212121212
$ws.Range("M18").Value = 0.03535353535353535
$ws.Range("O18").Value = 0.04545454545454546
$ws.Range("S18").Value = 0.07575757575757576
$ws.Range("F19").Value = 0.01776461880088823
$ws.Range("H19").Value = 0.2131754256106588
$ws.Range("I19").Value = 0.09474463360473723
$ws.Range("J19").Value = 0.3671354552183568
$ws.Range("K19").Value = 0.09474463360473723
$ws.Range("M19").Value = 0.01406365655070318
$ws.Range("O19").Value = 0.07920059215396003
$ws.Range("S19").Value = 0.1191709844559585
